$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2: refreshed counts/probabilities
$ws.Range("B2").Value = 0.2244897959183673
$ws.Range("C2").Value = 0.5131195335276968
$ws.Range("J2").Value = 0.01749271137026239
$ws.Range("P2").Value = 0.1749271137026239
$ws.Range("S2").Value = 0.06997084548104957

# Row 3: refreshed counts/probabilities
$ws.Range("C3").Value = 0.02702702702702703
$ws.Range("J3").Value = 0.04864864864864865
$ws.Range("P3").Value = 0.7513513513513513
$ws.Range("S3").Value = 0.172972972972973

# Row 4: refreshed counts/probabilities
$ws.Range("J4").Value = 0.0975609756097561
$ws.Range("P4").Value = 0.6341463414634146
$ws.Range("S4").Value = 0.2682926829268293

# Row 6: refreshed counts/probabilities
$ws.Range("B6").Value = 0.07870370370370371
$ws.Range("F6").Value = 0.05555555555555555
$ws.Range("J6").Value = 0.2824074074074074
$ws.Range("O6").Value = 0.03240740740740741
$ws.Range("Q6").Value = 0.1435185185185185
$ws.Range("R6").Value = 0.04166666666666666
$ws.Range("S6").Value = 0.3657407407407408

# Row 7: refreshed counts/probabilities
$ws.Range("B7").Value = 0.1517857142857143
$ws.Range("D7").Value = 0.02232142857142857
$ws.Range("F7").Value = 0.0625
$ws.Range("J7").Value = 0.1428571428571428
$ws.Range("O7").Value = 0.02678571428571428
$ws.Range("Q7").Value = 0.1785714285714286
$ws.Range("R7").Value = 0.05357142857142857
$ws.Range("S7").Value = 0.3616071428571428

# Row 8: refreshed counts/probabilities
$ws.Range("B8").Value = 0.1129032258064516
$ws.Range("D8").Value = 0.01290322580645161
$ws.Range("F8").Value = 0.06774193548387097
$ws.Range("J8").Value = 0.1258064516129032
$ws.Range("O8").Value = 0.03870967741935484
$ws.Range("Q8").Value = 0.1741935483870968
$ws.Range("R8").Value = 0.07741935483870968
$ws.Range("S8").Value = 0.3903225806451613

# Row 9: refreshed counts/probabilities
$ws.Range("B9").Value = 0.1012658227848101
$ws.Range("D9").Value = 0.02531645569620253
$ws.Range("F9").Value = 0.06962025316455696
$ws.Range("J9").Value = 0.1329113924050633
$ws.Range("O9").Value = 0.0379746835443038
$ws.Range("Q9").Value = 0.2025316455696203
$ws.Range("R9").Value = 0.08860759493670886
$ws.Range("S9").Value = 0.3417721518987342

# Row 10: refreshed counts/probabilities
$ws.Range("B10").Value = 0.1332788225674571
$ws.Range("D10").Value = 0.02534750613246116
$ws.Range("E10").Value = 0.001635322976287817
$ws.Range("F10").Value = 0.06950122649223221
$ws.Range("J10").Value = 0.1275551921504497
$ws.Range("O10").Value = 0.02371218315617334
$ws.Range("Q10").Value = 0.2134096484055601
$ws.Range("R10").Value = 0.05314799672935405
$ws.Range("S10").Value = 0.3524121013900245

# Row 11: refreshed counts/probabilities
$ws.Range("G11").Value = 0.1472222222222222
$ws.Range("J11").Value = 0.08888888888888889
$ws.Range("K11").Value = 0.2027777777777778
$ws.Range("L11").Value = 0.5472222222222223
$ws.Range("S11").Value = 0.01388888888888889

# Row 12: refreshed counts/probabilities
$ws.Range("G12").Value = 0.7560975609756098
$ws.Range("J12").Value = 0.1804878048780488
$ws.Range("K12").Value = 0.00975609756097561
$ws.Range("L12").Value = 0.02439024390243903
$ws.Range("S12").Value = 0.02926829268292683

# Row 13: refreshed counts/probabilities
$ws.Range("G13").Value = 0.5925925925925926
$ws.Range("J13").Value = 0.3333333333333333
$ws.Range("S13").Value = 0.07407407407407407

# Row 14: refreshed counts/probabilities
$ws.Range("G14").Value = 0.8
$ws.Range("J14").Value = 0.2

# Row 15: refreshed counts/probabilities
$ws.Range("F15").Value = 0.01388888888888889
$ws.Range("H15").Value = 0.1342592592592593
$ws.Range("I15").Value = 0.08333333333333333
$ws.Range("J15").Value = 0.3611111111111111
$ws.Range("K15").Value = 0.05092592592592592
$ws.Range("M15").Value = 0.004629629629629629
$ws.Range("O15").Value = 0.05092592592592592
$ws.Range("S15").Value = 0.3009259259259259

# Row 16: refreshed counts/probabilities
$ws.Range("F16").Value = 0.03603603603603604
$ws.Range("H16").Value = 0.1216216216216216
$ws.Range("I16").Value = 0.08558558558558559
$ws.Range("J16").Value = 0.454954954954955
$ws.Range("K16").Value = 0.1396396396396396
$ws.Range("M16").Value = 0.004504504504504504
$ws.Range("N16").Value = 0.004504504504504504
$ws.Range("O16").Value = 0.04054054054054054
$ws.Range("S16").Value = 0.1126126126126126

# Row 17: refreshed counts/probabilities
$ws.Range("F17").Value = 0.03125
$ws.Range("H17").Value = 0.1682692307692308
$ws.Range("I17").Value = 0.07211538461538461
$ws.Range("J17").Value = 0.4134615384615384
$ws.Range("K17").Value = 0.1201923076923077
$ws.Range("M17").Value = 0.02163461538461538
$ws.Range("O17").Value = 0.06490384615384616
$ws.Range("S17").Value = 0.1081730769230769

# Row 18: refreshed counts/probabilities
$ws.Range("F18").Value = 0.02479338842975207
$ws.Range("H18").Value = 0.140495867768595
$ws.Range("I18").Value = 0.04958677685950413
$ws.Range("J18").Value = 0.4380165289256198
$ws.Range("K18").Value = 0.140495867768595
$ws.Range("M18").Value = 0.02479338842975207
$ws.Range("O18").Value = 0.09917355371900827
$ws.Range("S18").Value = 0.08264462809917356

# Row 19: refreshed counts/probabilities
$ws.Range("F19").Value = 0.01818181818181818
$ws.Range("H19").Value = 0.1545454545454545
$ws.Range("I19").Value = 0.07818181818181819
$ws.Range("J19").Value = 0.3881818181818182
$ws.Range("K19").Value = 0.1581818181818182
$ws.Range("M19").Value = 0.01181818181818182
$ws.Range("N19").Value = 0.003636363636363636
$ws.Range("O19").Value = 0.07363636363636364
$ws.Range("S19").Value = 0.1136363636363636
